$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 254, pushing existing rows
# (254-271) down to (256-273). Formatting is inherited from the row above.
$ws.Rows(254).Insert()
$ws.Rows(254).Insert()

# --- New row 254: Ají / Americana (o) / Primera, Comercializadora del Agro de Limarí ---
$ws.Range("A254").Value = 2
$ws.Range("B254").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C254").Value = "Coquimbo"
$ws.Range("D254").Value = 44714
$ws.Range("E254").Value = 4
$ws.Range("F254").Value = 100112021
$ws.Range("G254").Value = "Ají"
$ws.Range("H254").Value = "Americana (o)"
$ws.Range("I254").Value = "Primera"
$ws.Range("J254").Value = 140
$ws.Range("K254").Value = 20000
$ws.Range("L254").Value = 23000
$ws.Range("M254").Value = 21500
$ws.Range("N254").Value = "`$/caja 25 kilos"
$ws.Range("O254").Value = "Provincia de Limarí"
$ws.Range("P254").Value = 860
$ws.Range("Q254").Value = 25
$ws.Range("R254").Value = "Hortaliza"

# --- New row 255: Ají / Inferno / Primera, Comercializadora del Agro de Limarí ---
$ws.Range("A255").Value = 2
$ws.Range("B255").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C255").Value = "Coquimbo"
$ws.Range("D255").Value = 44714
$ws.Range("E255").Value = 4
$ws.Range("F255").Value = 100112021
$ws.Range("G255").Value = "Ají"
$ws.Range("H255").Value = "Inferno"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 100
$ws.Range("K255").Value = 30000
$ws.Range("L255").Value = 33000
$ws.Range("M255").Value = 31500
$ws.Range("N255").Value = "`$/caja 25 kilos"
$ws.Range("O255").Value = "Provincia de Limarí"
$ws.Range("P255").Value = 1260
$ws.Range("Q255").Value = 25
$ws.Range("R255").Value = "Hortaliza"

# --- The former row 269 (now shifted to row 271) only changes its date ---
$ws.Range("D271").Value = 44714
